# Actualización automática 2025-11-10 16:30:09
# Updates the "CUMPLIMIENTO MENSUAL" sheet: new PRESUPUESTO / VENTA figures
# (and their derived POR CUMPLIR / CUMPLIMIENTO columns), plus a small
# column-width tweak on columns D:F.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# --- Column widths (cols D, E, F) ---
# ColumnWidth is expressed in "characters" and Excel stores the sheet's
# <col width="..."> a bit higher (~0.83 chars, the default cell-padding
# offset for this font/theme) than what was asked for, so back that
# padding out to land on the exact stored widths from the target file
# (12, 22, 26).
$ws.Columns.Item(4).ColumnWidth = 12 - 0.83
$ws.Columns.Item(5).ColumnWidth = 22 - 0.83
$ws.Columns.Item(6).ColumnWidth = 26 - 0.83

# --- Data rows: C = PRESUPUESTO, D = VENTA, E = POR CUMPLIR, F = CUMPLIMIENTO ---
function Set-Row($row, $c, $d, $e, $f) {
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
    $ws.Cells.Item($row, 6).Value = $f
}

Set-Row 2  129.6                0     129.6                0
Set-Row 3  2344.03              0     2344.03              0
Set-Row 4  207.39               0     207.39               0
Set-Row 5  86.41                0     86.41                0
Set-Row 6  855.91               0     855.91               0
Set-Row 7  383                  0     383                  0
Set-Row 8  415                  0     415                  0
Set-Row 10 388.107983534392     0     388.107983534392     0
Set-Row 11 902.88               443.44   459.44            0.4911394648236754
Set-Row 12 34701                317.92   34383.08          0.009161695628368058
Set-Row 13 364.412605947529     0     364.412605947529     0
Set-Row 14 40777.74058948192    761.36   40016.38058948192 0.01867097070592437
